$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/cc2ce6eb366f1034d22c26ed2ee6b05b84f3e7a6/e2e/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3087839134feca2713bf27c7a424e7afc32d48ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7f69f789c581b2b476ae44bb0ae79d3d67f47e62/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"

$mdName = "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
$zhXlfName = "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
$deXlfName = "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"

$status = "Handed back: in sync with en-US"

# ---------- zh-cn sheet ----------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("C2").Value = $status
$ws.Range("C3").Value = $status

$ws.Range("H2").Value = "2016-03-22 21:11:29"
$ws.Range("H3").Value = "2016-03-22 21:11:29"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("F2").Value = $mdName
$ws.Range("G2").Value = $zhXlfName
$ws.Range("F3").Value = $mdName
$ws.Range("G3").Value = $zhXlfName

$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("G2"), $zhXlfUrl, "", "", $zhXlfName)
$ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("G3"), $zhXlfUrl, "", "", $zhXlfName)

$ws.Range("F2:G3").Font.Underline = 2
$ws.Range("F2:G3").Font.Color = 13011546

# ---------- de-de sheet ----------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("C2").Value = $status
$ws2.Range("C3").Value = $status

$ws2.Range("H2").Value = "2016-03-22 21:11:38"
$ws2.Range("H3").Value = "2016-03-22 21:11:38"
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("F2").Value = $mdName
$ws2.Range("G2").Value = $deXlfName
$ws2.Range("F3").Value = $mdName
$ws2.Range("G3").Value = $deXlfName

$ws2.Hyperlinks.Add($ws2.Range("F2"), $mdUrl, "", "", $mdName)
$ws2.Hyperlinks.Add($ws2.Range("G2"), $deXlfUrl, "", "", $deXlfName)
$ws2.Hyperlinks.Add($ws2.Range("F3"), $mdUrl, "", "", $mdName)
$ws2.Hyperlinks.Add($ws2.Range("G3"), $deXlfUrl, "", "", $deXlfName)

$ws2.Range("F2:G3").Font.Underline = 2
$ws2.Range("F2:G3").Font.Color = 13011546

Write-Host "Edits applied"
